# Generated test cases update for the "Testcases" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testcases")

# --- Header block ---------------------------------------------------------
$ws.Range("B2").Value = "Component: Multi-Functional Tool Application"
$ws.Range("E3").Value = "MFP: Any"

# --- TC001 (row 6) ---------------------------------------------------------
$ws.Range("B6").Value = "TC001"
$ws.Range("C6").Value = "Application is installed on a desktop computer"
$ws.Range("D6").Value = "Verify application installation on desktop"
$ws.Range("E6").Value = "1. Copy MultiFunctionalTool_For_Desktop.zip from specified tec-share location`n2. Extract contents to preferred location`n3. Double-click MultiFunctionalToolApplication"
$ws.Range("F6").Value = "Application launches successfully with all features accessible"
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""

# --- TC002 (row 7) ---------------------------------------------------------
$ws.Range("B7").Value = "TC002"
$ws.Range("C7").Value = "Application is installed on a laptop"
$ws.Range("D7").Value = "Verify application installation on laptop"
$ws.Range("E7").Value = "1. Copy MultiFunctionalTool_For_Laptop.zip from specified tec-share location`n2. Extract contents to preferred location`n3. Double-click MultiFunctionalToolApplication"
$ws.Range("F7").Value = "Application launches successfully with all features accessible"
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = ""

# --- TC003 (row 8) ---------------------------------------------------------
$ws.Range("B8").Value = "TC003"
$ws.Range("C8").Value = "Application is installed and running"
$ws.Range("D8").Value = "Verify Network Packet Capture start functionality"
$ws.Range("E8").Value = "1. Navigate to Network Packet Capture section`n2. Click Start button"
$ws.Range("F8").Value = "Packet capture begins successfully"
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = ""

# --- TC004 (row 9) ---------------------------------------------------------
$ws.Range("B9").Value = "TC004"
$ws.Range("C9").Value = "Network Packet Capture is running"
$ws.Range("D9").Value = "Verify Network Packet Capture stop functionality"
$ws.Range("E9").Value = "1. Navigate to Network Packet Capture section`n2. Click Stop button"
$ws.Range("F9").Value = "1. Packet capture stops`n2. .pcap file is generated`n3. File is copied to MFP's Shared Folder`n4. Shared Folder opens automatically"
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = ""

# --- TC005 (row 10) ---------------------------------------------------------
$ws.Range("B10").Value = "TC005"
$ws.Range("C10").Value = "Application is installed and running"
$ws.Range("D10").Value = "Verify Memory Leak Check functionality"
$ws.Range("E10").Value = "1. Navigate to Memory Leak Check section`n2. Select a protocol`n3. Run the memory leak check"
$ws.Range("F10").Value = "Memory Leak Comparison Table is displayed with accurate information"
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = ""

# --- TC006 (row 11) ---------------------------------------------------------
$ws.Range("B11").Value = "TC006"
$ws.Range("C11").Value = "Application is installed and running"
$ws.Range("D11").Value = "Verify Debug Log Collection functionality"
$ws.Range("E11").Value = "1. Navigate to Debug Log Collection section`n2. Click Run button"
$ws.Range("F11").Value = "1. Script executes successfully`n2. Logs are collected`n3. Logs are copied to MFP's Shared Folder`n4. Shared Folder opens automatically"
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = ""

# --- TC007 (row 12) ---------------------------------------------------------
$ws.Range("B12").Value = "TC007"
$ws.Range("C12").Value = "Debug Log Collection has been run once with empty folder result"
$ws.Range("D12").Value = "Verify Debug Log Collection retry functionality"
$ws.Range("E12").Value = "1. Navigate to Debug Log Collection section`n2. Click Run button again"
$ws.Range("F12").Value = "1. Script executes successfully`n2. Logs are collected`n3. Logs are copied to MFP's Shared Folder`n4. Shared Folder opens with logs visible"
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = ""

# --- TC008 (row 13) ---------------------------------------------------------
$ws.Range("B13").Value = "TC008"
$ws.Range("C13").Value = "Application is installed and running"
$ws.Range("D13").Value = "Verify Diagnostic Code Details - ECC selection"
$ws.Range("E13").Value = "1. Navigate to Diagnostic Code Details section`n2. Select ECC option`n3. Choose specific diagnostic code"
$ws.Range("F13").Value = "Relevant job-specific details for ECC are displayed correctly"
$ws.Range("G13").Value = ""
$ws.Range("H13").Value = ""

# --- TC009 (row 14) ---------------------------------------------------------
$ws.Range("B14").Value = "TC009"
$ws.Range("C14").Value = "Application is installed and running"
$ws.Range("D14").Value = "Verify Diagnostic Code Details - Network Protocols selection"
$ws.Range("E14").Value = "1. Navigate to Diagnostic Code Details section`n2. Select Network Protocols option`n3. Choose specific diagnostic code"
$ws.Range("F14").Value = "Relevant job-specific details for Network Protocols are displayed correctly"
$ws.Range("G14").Value = ""
$ws.Range("H14").Value = ""

# --- TC010 (row 15) ---------------------------------------------------------
$ws.Range("B15").Value = "TC010"
$ws.Range("C15").Value = "Application is installed and running"
$ws.Range("D15").Value = "Verify Diagnostic Code Details - High Security Mode selection"
$ws.Range("E15").Value = "1. Navigate to Diagnostic Code Details section`n2. Select High Security Mode option`n3. Choose specific diagnostic code"
$ws.Range("F15").Value = "Relevant job-specific details for High Security Mode are displayed correctly"
$ws.Range("G15").Value = ""
$ws.Range("H15").Value = ""

# --- TC011 (row 16) ---------------------------------------------------------
$ws.Range("B16").Value = "TC011"
$ws.Range("C16").Value = "Application is installed and running"
$ws.Range("D16").Value = "Verify 08 Diagnostic Code Value - Get functionality"
$ws.Range("E16").Value = "1. Navigate to 08 Diagnostic Code Value section`n2. Select a diagnostic code`n3. Click Get button"
$ws.Range("F16").Value = "Current value of the selected diagnostic code is displayed"
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = ""

# --- TC012 (row 17) ---------------------------------------------------------
$ws.Range("B17").Value = "TC012"
$ws.Range("C17").Value = "Application is installed and running"
$ws.Range("D17").Value = "Verify 08 Diagnostic Code Value - Set functionality"
$ws.Range("E17").Value = "1. Navigate to 08 Diagnostic Code Value section`n2. Select a diagnostic code`n3. Enter a new value`n4. Click Set button"
$ws.Range("F17").Value = "Diagnostic code value is updated successfully"
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = ""

# --- TC013 (row 18) ---------------------------------------------------------
$ws.Range("B18").Value = "TC013"
$ws.Range("C18").Value = "Application is installed and running"
$ws.Range("D18").Value = "Verify Protocol Configuration - Get functionality"
$ws.Range("E18").Value = "1. Navigate to Protocol Configuration section`n2. Select a protocol`n3. Click Get button"
$ws.Range("F18").Value = "Current value of the selected protocol is displayed"
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = ""

# --- TC014 (row 19) ---------------------------------------------------------
$ws.Range("B19").Value = "TC014"
$ws.Range("C19").Value = "Application is installed and running"
$ws.Range("D19").Value = "Verify Protocol Configuration - Set functionality placeholder"
$ws.Range("E19").Value = "1. Navigate to Protocol Configuration section`n2. Select a protocol`n3. Attempt to set a value"
$ws.Range("F19").Value = "System indicates that the Set protocol values operation is not yet implemented"
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = ""

# --- TC015 (row 20) ---------------------------------------------------------
$ws.Range("B20").Value = "TC015"
$ws.Range("C20").Value = "Application is running with multiple features"
$ws.Range("D20").Value = "Verify simultaneous operation of multiple features"
$ws.Range("E20").Value = "1. Start Network Packet Capture`n2. While capture is running, perform Memory Leak Check`n3. While both are running, check Diagnostic Code Details"
$ws.Range("F20").Value = "All features operate correctly without interference"
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = ""

# --- TC016 (row 21) ---------------------------------------------------------
$ws.Range("B21").Value = "TC016"
$ws.Range("C21").Value = "Application is installed and running"
$ws.Range("D21").Value = "Verify application performance - startup time"
$ws.Range("E21").Value = "1. Close the application if running`n2. Time the startup of the application"
$ws.Range("F21").Value = "Application starts in under 5 seconds"
$ws.Range("G21").Value = ""
$ws.Range("H21").Value = "Non-functional test"

# --- TC017 (row 22) ---------------------------------------------------------
$ws.Range("B22").Value = "TC017"
$ws.Range("C22").Value = "Application is installed and running"
$ws.Range("D22").Value = "Verify application performance - response time"
$ws.Range("E22").Value = "1. Click various buttons and features`n2. Measure response time"
$ws.Range("F22").Value = "UI responds within 1 second to user interactions"
$ws.Range("G22").Value = ""
$ws.Range("H22").Value = "Non-functional test"

# --- TC018 (row 23, brand-new row) ------------------------------------------
$ws.Range("B23").Value = "TC018"
$ws.Range("C23").Value = "Application is installed and running"
$ws.Range("D23").Value = "Verify application usability - intuitive interface"
$ws.Range("E23").Value = "1. Observe the layout of the application`n2. Attempt to use each feature without prior training"
$ws.Range("F23").Value = "User can navigate and use all features without confusion"
$ws.Range("G23").Value = ""
$ws.Range("H23").Value = "Non-functional test"

# --- TC019 (row 24, reclaims the old "Test Summary" label) ------------------
$ws.Range("B24").Value = "TC019"
$ws.Range("C24").Value = "Application is installed and running"
$ws.Range("D24").Value = "Verify application compatibility with different screen resolutions"
$ws.Range("E24").Value = "1. Run application on monitors with different resolutions`n2. Check if all UI elements are visible and properly sized"
$ws.Range("F24").Value = "Application displays correctly on all tested resolutions"
$ws.Range("G24").Value = ""
$ws.Range("H24").Value = "Non-functional test"

# --- TC020 (row 25, reclaims the old "Test Case Count:" label) --------------
$ws.Range("B25").Value = "TC020"
$ws.Range("C25").Value = "Application is installed and running"
$ws.Range("D25").Value = "Verify application stability during extended use"
$ws.Range("E25").Value = "1. Run the application continuously for 8 hours`n2. Periodically perform operations with all features"
$ws.Range("F25").Value = "Application remains stable without crashes or memory issues"
$ws.Range("G25").Value = ""
$ws.Range("H25").Value = "Non-functional test"

# Rows 22-25 had no explicit row height before the edit and none afterwards
# either; writing multi-line text into them would otherwise leave a stamped
# auto-height behind, so re-fit (which clears the explicit height) to match.
$ws.Rows("22:25").AutoFit()
